$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 105-134 entirely; this shifts the former rows 135-138
# (the SOLEVUL / HOLRM-3 / HOLRA-3 records) up into rows 105-108,
# matching the dimension shrinking from A1:Q138 to A1:Q108.
$ws.Range("A105:A134").EntireRow.Delete()
